$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.95"
$ws.Range("D3").Value = "'22.79"
$ws.Range("D4").Value = "'5.287"
$ws.Range("D5").Value = "'0.05726"
$ws.Range("D6").Value = "'3.424"
$ws.Range("D7").Value = "'0.8103"
$ws.Range("D8").Value = "'0.8744"
$ws.Range("D10").Value = "'0.07418"
$ws.Range("D12").Value = "'0.03113"
$ws.Range("D13").Value = "'0.09388"
$ws.Range("D14").Value = "'3.878"
$ws.Range("D15").Value = "'0.001586"
$ws.Range("D16").Value = "'0.04815"
$ws.Range("D17").Value = "'0.0005840"
$ws.Range("E17").Value = "16OneONEWorstin24h"
$ws.Range("D18").Value = "'0.006145"
$ws.Range("D20").Value = "'0.0009968"
$ws.Range("D22").Value = "'3.736"
$ws.Range("D23").Value = "'6.333"
$ws.Range("D24").Value = "'2.198"
$ws.Range("D40").Value = "'0.03939"
$ws.Range("D41").Value = "'0.006763"
$ws.Range("D43").Value = "'0.003200"
$ws.Range("D44").Value = "'0.007252"
$ws.Range("D45").Value = "'0.00005615"
$ws.Range("D47").Value = "'0.6000"
$ws.Range("D48").Value = "'0.1796"
$ws.Range("E48").Value = "47BOLOBOLO"
